$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 101, shifting existing rows 101-118 down to 102-119.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly record.
$ws.Cells.Item(101, 1).Value = 10
$ws.Cells.Item(101, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(101, 3).Value = "La Araucanía"
$ws.Cells.Item(101, 4).Value = 44463
$ws.Cells.Item(101, 5).Value = 9
$ws.Cells.Item(101, 6).Value = 100112013
$ws.Cells.Item(101, 7).Value = "Alcachofa"
$ws.Cells.Item(101, 8).Value = "Madrigal"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 50
$ws.Cells.Item(101, 11).Value = 12000
$ws.Cells.Item(101, 12).Value = 12000
$ws.Cells.Item(101, 13).Value = 12000
$ws.Cells.Item(101, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(101, 15).Value = "Región Metropolitana"
$ws.Cells.Item(101, 16).Value = 300
$ws.Cells.Item(101, 17).Value = 40
$ws.Cells.Item(101, 18).Value = "Hortaliza"
